# Update the REST API mapping sheet: switch the domain/host/organization
# URI placeholders from "[name]" to "[id]" (and give /organization/
# the missing "/[id]" on the UPDATE row), matching the new
# ORGANIZATION* REST endpoints that were added to the mapping table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = "/contact/[id]"

$ws.Range("C7").Value  = "/domain/[id]"
$ws.Range("C8").Value  = "/domain/[id]"
$ws.Range("C9").Value  = "/domain/[id]"
$ws.Range("C10").Value = "/domain/[id]"

$ws.Range("C11").Value = "/domain/transfer/[id]"
$ws.Range("C13").Value = "/domain/renew/[id]"
$ws.Range("C12").Value = "/domain/querytransfer/[id]"

$ws.Range("C15").Value = "/host/[id]"
$ws.Range("C17").Value = "/host/[id]"

$ws.Range("C21").Value = "/organization/[id]"

$ws.Range("A6").Select()
